$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set text format on D and E columns for the data rows so Excel
# keeps the values as literal strings (matching the source price
# and volume formatting) instead of coercing them to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "26.666.67"
$ws.Range("E2").Value = "  -0.23%  "
$ws.Range("D3").Value = "1.598.95"
$ws.Range("E3").Value = "  -0.04%  "
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("D5").Value = "211.82"
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("E6").Value = "  +0.59%  "
$ws.Range("E7").Value = "  +0.35%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E9").Value = "  +0.23%  "
$ws.Range("D10").Value = "19.58"
$ws.Range("E10").Value = "  -0.94%  "
$ws.Range("E11").Value = "  +0.19%  "
$ws.Range("D12").Value = "1.823.85"
$ws.Range("E12").Value = "  +0.02%  "
$ws.Range("D13").Value = "1.588.32"
$ws.Range("E13").Value = "  -0.64%  "
$ws.Range("E14").Value = "  -0.10%  "
$ws.Range("E15").Value = "  +0.00%  "
$ws.Range("D16").Value = "65.14"
$ws.Range("E16").Value = "  +0.00%  "
$ws.Range("D17").Value = "26.670.95"
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("E18").Value = "  +1.00%  "
$ws.Range("E19").Value = "  +0.31%  "
$ws.Range("D20").Value = "208.84"
$ws.Range("E20").Value = "  -0.74%  "
$ws.Range("D21").Value = "7.04"
$ws.Range("E21").Value = "  +4.66%  "
$ws.Range("D22").Value = "4.29"
$ws.Range("E22").Value = "  +0.54%  "
$ws.Range("E23").Value = "  +0.64%  "
$ws.Range("D24").Value = "8.96"
$ws.Range("E24").Value = "  +0.44%  "
$ws.Range("D25").Value = "145.31"
$ws.Range("E25").Value = "  -1.05%  "
$ws.Range("E26").Value = "  +0.29%  "
$ws.Range("D27").Value = "7.11"
$ws.Range("E27").Value = "  -0.87%  "
$ws.Range("E28").Value = "  -0.57%  "
$ws.Range("D29").Value = "15.29"
$ws.Range("E29").Value = "  -0.18%  "
$ws.Range("E30").Value = "  +2.01%  "
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("E32").Value = "  +0.17%  "
$ws.Range("E33").Value = "  +1.16%  "
$ws.Range("D34").Value = "1.276.56"
$ws.Range("E34").Value = "  -1.97%  "
$ws.Range("E35").Value = "  -7.55%  "
$ws.Range("D36").Value = "2.46"
$ws.Range("E36").Value = "  +0.31%  "
$ws.Range("E37").Value = "  +0.80%  "
$ws.Range("E38").Value = "  -1.10%  "
$ws.Range("E39").Value = "  -1.21%  "
$ws.Range("E40").Value = "  +18.64%  "
$ws.Range("E41").Value = "  +2.40%  "
$ws.Range("E42").Value = "  +0.47%  "
$ws.Range("D43").Value = "0.785"
$ws.Range("E43").Value = "  -0.80%  "
$ws.Range("D44").Value = "63.97"
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("D45").Value = "1.735.97"
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("D46").Value = "91.24"
$ws.Range("E46").Value = "  +1.37%  "
$ws.Range("E47").Value = "  -2.51%  "
$ws.Range("E48").Value = "  +3.19%  "
$ws.Range("E49").Value = "  +0.72%  "
$ws.Range("E50").Value = "  +0.07%  "
$ws.Range("D51").Value = "7.40"
$ws.Range("E51").Value = "  -1.85%  "
